$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "31/12/2025"
$ws.Range("C7").Value = 1.52847110399664

$ws.Range("B13").Value = "31/12/2025"
$ws.Range("C13").Value = 1.52290951123554

$ws.Range("B19").Value = "31/12/2025"
$ws.Range("C19").Value = 1.50292233927388
